$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: ATC05_PagarCuenta -> ATC05_descargarCartola6meses
#   - clear old D6:F6 (Directv / 126554850 / Directv)
#   - add rut (B6) and apellido (C6) values, mirroring row 2's "rut"/"Rojas651" pair
$ws.Range("D6:F6").Clear()
$ws.Range("A6").Value = "ATC05_descargarCartola6meses"
$ws.Range("B6").Value = "175553878"
$ws.Range("C6").Value = "Rojas651"

# Row 7: ATC06 -> ATC06_descargarCartolaLuz, with the same new rut/apellido pair
$ws.Range("A7").Value = "ATC06_descargarCartolaLuz"
$ws.Range("B7").Value = "175553878"
$ws.Range("C7").Value = "Rojas651"

# Column A got wider (and no longer auto "best fit") to accommodate the longer
# new test-case names.
$ws.Columns("A").ColumnWidth = 28.3

# Selection moved to D5
[void]$ws.Range("D5").Select()
